# Updated test data as per new implementation:
# The "Loading Details Name" value on the "Add Panels" sheet changes from
# "Main Processor 24V (A)" to "24V Rail(A)" (cells K8 and K9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

$ws.Range("K8").Value = "24V Rail(A)"
$ws.Range("K9").Value = "24V Rail(A)"

# Move the active selection to H8, matching the saved cursor position.
$ws.Range("H8").Select()
